$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7058433170332705
$ws.Range("C2").Value = 0.07589299558478046
$ws.Range("E2").Value = 0.1676849289299227
$ws.Range("F2").Value = 2.400560874658339
$ws.Range("G2").Value = 1.322996667511106
$ws.Range("H2").Value = 1.259991072425294
$ws.Range("I2").Value = 1.285890466443668
$ws.Range("J2").Value = 0.08798896417971491
$ws.Range("K2").Value = 0.4149592270594553
$ws.Range("L2").Value = 0.4087342123875288
$ws.Range("M2").Value = 0.257632795613997
$ws.Range("N2").Value = 2.451327081456153

$ws.Range("B3").Value = 0.6727567676213937
$ws.Range("C3").Value = 0.07339633752365415
$ws.Range("E3").Value = 0.1675891279756989
$ws.Range("F3").Value = 2.400603688313481
$ws.Range("G3").Value = 1.327016002299288
$ws.Range("H3").Value = 1.266380605452071
$ws.Range("I3").Value = 1.292827113541033
$ws.Range("J3").Value = 0.08709307933932564
$ws.Range("K3").Value = 0.3831244869144257
$ws.Range("L3").Value = 0.4046225910701509
$ws.Range("M3").Value = 0.2505718341667276
$ws.Range("N3").Value = 2.472917094709249

$ws.Range("B4").Value = 0.6527413314738055
$ws.Range("C4").Value = 0.07183861730952401
$ws.Range("E4").Value = 0.167575824062915
$ws.Range("F4").Value = 2.401675821894329
$ws.Range("G4").Value = 1.330137978200739
$ws.Range("H4").Value = 1.270763307915431
$ws.Range("I4").Value = 1.297598862277525
$ws.Range("J4").Value = 0.08653640024180831
$ws.Range("K4").Value = 0.3637292863464978
$ws.Range("L4").Value = 0.4022535390827855
$ws.Range("M4").Value = 0.2463467320889201
$ws.Range("N4").Value = 2.486858812436656

$ws.Range("B5").Value = 0.6446608164006875
$ws.Range("C5").Value = 0.07119758485291783
$ws.Range("E5").Value = 0.1675818893520837
$ws.Range("F5").Value = 2.402376052257068
$ws.Range("G5").Value = 1.331574706989286
$ws.Range("H5").Value = 1.272664952329492
$ws.Range("I5").Value = 1.29967237733641
$ws.Range("J5").Value = 0.08630789975104491
$ws.Range("K5").Value = 0.3558640228125256
$ws.Range("L5").Value = 0.4013273832228137
$ws.Range("M5").Value = 0.2446528664425784
$ws.Range("N5").Value = 2.492712551790685

$ws.Range("B6").Value = 0.6433236554084658
$ws.Range("C6").Value = 0.07109076384069368
$ws.Range("E6").Value = 0.1675835914598078
$ws.Range("F6").Value = 2.402508239333571
$ws.Range("G6").Value = 1.331823210689052
$ws.Range("H6").Value = 1.272987707668321
$ws.Range("I6").Value = 1.300024476088613
$ws.Range("J6").Value = 0.08626985812938415
$ws.Range("K6").Value = 0.3545603337672674
$ws.Range("L6").Value = 0.4011759708978815
$ws.Range("M6").Value = 0.2443732906455942
$ws.Range("N6").Value = 2.493694972957595

$ws.Range("B7").Value = 0.6526320466358868
$ws.Range("C7").Value = 0.0718299974512675
$ws.Range("E7").Value = 0.1675758592977079
$ws.Range("F7").Value = 2.401684198792751
$ws.Range("G7").Value = 1.330156688351977
$ws.Range("H7").Value = 1.270788485717247
$ws.Range("I7").Value = 1.297626304042996
$ws.Range("J7").Value = 0.08653332526758106
$ws.Range("K7").Value = 0.3636230565498124
$ws.Range("L7").Value = 0.4022408894907059
$ws.Range("M7").Value = 0.2463237748714597
$ws.Range("N7").Value = 2.486937060021972

$ws.Range("B8").Value = 0.6943731558397417
$ws.Range("C8").Value = 0.0750372733232112
$ws.Range("E8").Value = 0.1676424736282982
$ws.Range("F8").Value = 2.400358770330939
$ws.Range("G8").Value = 1.324246789273161
$ws.Range("H8").Value = 1.262098879742055
$ws.Range("I8").Value = 1.288175879208602
$ws.Range("J8").Value = 0.08768143970787534
$ws.Range("K8").Value = 0.4039513745442207
$ws.Range("L8").Value = 0.4072843393693262
$ws.Range("M8").Value = 0.2551753573700246
$ws.Range("N8").Value = 2.458629094827444

$ws.Range("B9").Value = 0.7785873981084137
$ws.Range("C9").Value = 0.08113148788276447
$ws.Range("E9").Value = 0.1681325936292346
$ws.Range("F9").Value = 2.406045444771621
$ws.Range("G9").Value = 1.317847542341269
$ws.Range("H9").Value = 1.24869998420148
$ws.Range("I9").Value = 1.273707647701094
$ws.Range("J9").Value = 0.08988011221475745
$ws.Range("K9").Value = 0.4842255846492947
$ws.Range("L9").Value = 0.418402602967447
$ws.Range("M9").Value = 0.2734035403171617
$ws.Range("N9").Value = 2.408552857561897

$ws.Range("B10").Value = 0.8418797389525707
$ws.Range("C10").Value = 0.08549191904215547
$ws.Range("E10").Value = 0.168709753808475
$ws.Range("F10").Value = 2.415261617648952
$ws.Range("G10").Value = 1.316311904407002
$ws.Range("H10").Value = 1.241070101419069
$ws.Range("I10").Value = 1.265551443918888
$ws.Range("J10").Value = 0.09146290897181331
$ws.Range("K10").Value = 0.54391994117114
$ws.Range("L10").Value = 0.4273133865803231
$ws.Range("M10").Value = 0.2873207327667799
$ws.Range("N10").Value = 2.3750709855316

$ws.Range("B11").Value = 0.870977910827321
$ws.Range("C11").Value = 0.08745064425026783
$ws.Range("E11").Value = 0.1690190230540694
$ws.Range("F11").Value = 2.420545766548699
$ws.Range("G11").Value = 1.316301285316541
$ws.Range("H11").Value = 1.23807875342726
$ws.Range("I11").Value = 1.262377328860964
$ws.Range("J11").Value = 0.09217581533350483
$ws.Range("K11").Value = 0.5712306129101989
$ws.Range("L11").Value = 0.4315269996002655
$ws.Range("M11").Value = 0.2937649271723828
$ws.Range("N11").Value = 2.360556798609267

$ws.Range("B12").Value = 0.8820402059043886
$ws.Range("C12").Value = 0.08818881652337041
$ws.Range("E12").Value = 0.1691428141861167
$ws.Range("F12").Value = 2.422703439157829
$ws.Range("G12").Value = 1.316396208976855
$ws.Range("H12").Value = 1.237014867678496
$ws.Range("I12").Value = 1.261252409737857
$ws.Range("J12").Value = 0.09244474134396441
$ws.Range("K12").Value = 0.581594530325475
$ws.Range("L12").Value = 0.4331454659664473
$ws.Range("M12").Value = 0.2962213260387756
$ws.Range("N12").Value = 2.355163712337077

$ws.Range("B13").Value = 0.8796558195273008
$ws.Range("C13").Value = 0.08802999568658265
$ws.Range("E13").Value = 0.1691158572022289
$ws.Range("F13").Value = 2.422231782143939
$ws.Range("G13").Value = 1.316371364570614
$ws.Range("H13").Value = 1.237240932670886
$ws.Range("I13").Value = 1.261491255245005
$ws.Range("J13").Value = 0.09238686960550169
$ws.Range("K13").Value = 0.5793615048524146
$ws.Range("L13").Value = 0.4327958858917498
$ws.Range("M13").Value = 0.2956915820139443
$ws.Range("N13").Value = 2.356320623033454

$ws.Range("B14").Value = 0.8718871450299162
$ws.Range("C14").Value = 0.08751144525071197
$ws.Range("E14").Value = 0.1690290738260174
$ws.Range("F14").Value = 2.420720141189065
$ws.Range("G14").Value = 1.316307111633662
$ws.Range("H14").Value = 1.23798984705995
$ws.Range("I14").Value = 1.262283237205416
$ws.Range("J14").Value = 0.09219796087828414
$ws.Range("K14").Value = 0.5720828207770126
$ws.Range("L14").Value = 0.4316596944779718
$ws.Range("M14").Value = 0.2939666943574224
$ws.Range("N14").Value = 2.360111039244131

$ws.Range("B15").Value = 0.8671342487664617
$ws.Range("C15").Value = 0.08719335583172949
$ws.Range("E15").Value = 0.1689767849303685
$ws.Range("F15").Value = 2.419814612924952
$ws.Range("G15").Value = 1.316280640951973
$ws.Range("H15").Value = 1.238457545598578
$ws.Range("I15").Value = 1.26277838156151
$ws.Range("J15").Value = 0.09208211350598106
$ws.Range("K15").Value = 0.5676272647386895
$ws.Range("L15").Value = 0.4309667172748419
$ws.Range("M15").Value = 0.2929122458671927
$ws.Range("N15").Value = 2.362446209126592

$ws.Range("B16").Value = 0.8399842164927804
$ws.Range("C16").Value = 0.0853634132642469
$ws.Range("E16").Value = 0.1686904789509285
$ws.Range("F16").Value = 2.414938227541029
$ws.Range("G16").Value = 1.316326443279365
$ws.Range("H16").Value = 1.241275235249745
$ws.Range("I16").Value = 1.265769664619704
$ws.Range("J16").Value = 0.09141617457395768
$ws.Range("K16").Value = 0.5421382185421919
$ws.Range("L16").Value = 0.4270412246670787
$ws.Range("M16").Value = 0.2869018540203498
$ws.Range("N16").Value = 2.376033950181448

$ws.Range("B17").Value = 0.8234065571045051
$ws.Range("C17").Value = 0.08423445054853573
$ws.Range("E17").Value = 0.1685267760463951
$ws.Range("F17").Value = 2.412226111278386
$ws.Range("G17").Value = 1.316530755142196
$ws.Range("H17").Value = 1.243126555434401
$ws.Range("I17").Value = 1.267742011205321
$ws.Range("J17").Value = 0.09100581091073323
$ws.Range("K17").Value = 0.5265410417923988
$ws.Range("L17").Value = 0.4246739467161404
$ws.Range("M17").Value = 0.283243554778764
$ws.Range("N17").Value = 2.38455322387588

$ws.Range("B18").Value = 0.8139003784580154
$ws.Range("C18").Value = 0.08358276022826772
$ws.Range("E18").Value = 0.1684370180094987
$ws.Range("F18").Value = 2.410768944403202
$ws.Range("G18").Value = 1.316713021930155
$ws.Range("H18").Value = 1.244236525747056
$ws.Range("I18").Value = 1.2689269239944
$ws.Range("J18").Value = 0.09076911158516765
$ws.Range("K18").Value = 0.5175846094165877
$ws.Range("L18").Value = 0.423327425121812
$ws.Range("M18").Value = 0.2811500614014548
$ws.Range("N18").Value = 2.389520759090754

$ws.Range("B19").Value = 0.8106867227699297
$ws.Range("C19").Value = 0.08336170643620733
$ws.Range("E19").Value = 0.1684073843309974
$ws.Range("F19").Value = 2.410293232845717
$ws.Range("G19").Value = 1.316785854998258
$ws.Range("H19").Value = 1.244620097967399
$ws.Range("I19").Value = 1.269336785214584
$ws.Range("J19").Value = 0.09068885475174682
$ws.Range("K19").Value = 0.5145546435023505
$ws.Range("L19").Value = 0.4228741096898432
$ws.Range("M19").Value = 0.2804430767139507
$ws.Range("N19").Value = 2.391214268954361

$ws.Range("B20").Value = 0.8251682952993917
$ws.Range("C20").Value = 0.08435487267945518
$ws.Range("E20").Value = 0.1685437474550682
$ws.Range("F20").Value = 2.412504186478628
$ws.Range("G20").Value = 1.316502304071136
$ws.Range("H20").Value = 1.242924808024995
$ws.Range("I20").Value = 1.267526828349411
$ws.Range("J20").Value = 0.09104956415503906
$ws.Range("K20").Value = 0.5281998743606664
$ws.Range("L20").Value = 0.4249243884934657
$ws.Range("M20").Value = 0.2836318844585435
$ws.Range("N20").Value = 2.383639349067899

$ws.Range("B21").Value = 0.8741678184591706
$ws.Range("C21").Value = 0.08766385239010788
$ws.Range("E21").Value = 0.1690543832991871
$ws.Range("F21").Value = 2.421159896918681
$ws.Range("G21").Value = 1.316323298794529
$ws.Range("H21").Value = 1.237768004211517
$ws.Range("I21").Value = 1.262048522227744
$ws.Range("J21").Value = 0.09225347617458723
$ws.Range("K21").Value = 0.5742201538755864
$ws.Range("L21").Value = 0.4319928021725161
$ws.Range("M21").Value = 0.2944728993533232
$ws.Range("N21").Value = 2.358994903150823

$ws.Range("B22").Value = 0.9064448247261225
$ws.Range("C22").Value = 0.08980575304077831
$ws.Range("E22").Value = 0.1694270168219489
$ws.Range("F22").Value = 2.427730016376202
$ws.Range("G22").Value = 1.316783043762328
$ws.Range("H22").Value = 1.234799140826624
$ws.Range("I22").Value = 1.258917210996493
$ws.Range("J22").Value = 0.09303425673071075
$ws.Range("K22").Value = 0.6044248784457125
$ws.Range("L22").Value = 0.4367456087197041
$ws.Range("M22").Value = 0.301652030545668
$ws.Range("N22").Value = 2.343489421151656

$ws.Range("B23").Value = 0.8891950217559668
$ws.Range("C23").Value = 0.08866446865206967
$ws.Range("E23").Value = 0.1692245888699055
$ws.Range("F23").Value = 2.424139969274478
$ws.Range("G23").Value = 1.31648489172801
$ws.Range("H23").Value = 1.236346977316799
$ws.Range("I23").Value = 1.260547377380348
$ws.Range("J23").Value = 0.09261809690317335
$ws.Range("K23").Value = 0.5882924823526992
$ws.Range("L23").Value = 0.4341968100550986
$ws.Range("M23").Value = 0.2978118523024378
$ws.Range("N23").Value = 2.351709977001995

$ws.Range("B24").Value = 0.8243717369823287
$ws.Range("C24").Value = 0.08430043804177245
$ws.Range("E24").Value = 0.1685360611071509
$ws.Range("F24").Value = 2.412378150731925
$ws.Range("G24").Value = 1.316514964935564
$ws.Range("H24").Value = 1.243015875910544
$ws.Range("I24").Value = 1.267623953698703
$ws.Range("J24").Value = 0.09102978572861886
$ws.Range("K24").Value = 0.5274498831340679
$ws.Range("L24").Value = 0.4248111187204415
$ws.Range("M24").Value = 0.2834562903483544
$ws.Range("N24").Value = 2.384052294780346

$ws.Range("B25").Value = 0.7555544794047648
$ws.Range("C25").Value = 0.07950353271430544
$ws.Range("E25").Value = 0.1679617332793981
$ws.Range("F25").Value = 2.40362134793925
$ws.Range("G25").Value = 1.319022883914698
$ws.Range("H25").Value = 1.251935487460628
$ws.Range("I25").Value = 1.277187015631988
$ws.Range("J25").Value = 0.08929100363172182
$ws.Range("K25").Value = 0.4623827483016782
$ws.Range("L25").Value = 0.415264013537552
$ws.Range("M25").Value = 0.2683797444470528
$ws.Range("N25").Value = 2.421518371045835
